$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7, shifting existing rows 7-15 down to 8-16
$ws.Rows("7:7").Insert()

# Populate the newly inserted row 7 with the MONO buffer entry
$ws.Range("A7").Value = "#1Mosaique_corpusAudioMONO"
$ws.Range("B7").Value = "buffer audio mono pour l'analyse"

# Keep the active selection on A8, matching the saved sheet view state
$ws.Range("A8").Select()
